$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 24 de Junio de 2020 a las 20:13'

$data = @(
    @('Estados Unidos', 2441012, 16844, 1025285, 1291887, 0, 367, 123840),
    @('Brasil', 1157451, 5972, 613345, 491155, 0, 180, 52951),
    @('Rusia', 606881, 7176, 368822, 229546, 0, 154, 8513),
    @('India', 472947, 16832, 271432, 186609, 0, 423, 14906),
    @('Reino Unido', 306862, 652, 0, 0, 0, 154, 43081),
    @('España', 294166, 334, 0, 0, 0, 2, 28327),
    @('Peru', 260810, 0, 148437, 103969, 0, 0, 8404),
    @('Chile', 254416, 3649, 210570, 39115, 0, 226, 4731),
    @('Italia', 238833, 0, 184585, 19573, 0, 0, 34675),
    @('Iran', 212501, 2531, 172096, 30409, 0, 133, 9996),
    @('Alemania', 192913, 135, 176300, 7626, 0, 1, 8987),
    @('Turquia', 191657, 1492, 164234, 22398, 0, 24, 5025),
    @('Mexico', 191410, 6288, 144448, 23585, 0, 793, 23377),
    @('Pakistan', 188926, 3892, 77754, 107417, 0, 60, 3755),
    @('Arabia Saudita', 167267, 3123, 112797, 53083, 0, 41, 1387),
    @('Francia', 161267, 0, 74871, 56676, 0, 0, 29720),
    @('Banglades', 122660, 3462, 49666, 71412, 0, 37, 1582),
    @('Sudafrica', 106108, 0, 55045, 48961, 0, 0, 2102),
    @('Canada', 102179, 216, 64933, 28763, 0, 29, 8483),
    @('Catar', 90778, 1199, 73083, 17591, 0, 5, 104),
    @('China', 83430, 12, 78428, 368, 0, 0, 4634),
    @('Colombia', 73572, 0, 30459, 40709, 0, 0, 2404),
    @('Belgica', 60898, 88, 16771, 34405, 0, 9, 9722),
    @('Suecia', 60837, 0, 0, 0, 0, 0, 5161),
    @('Bielorrusia', 59945, 458, 40136, 19447, 0, 5, 362),
    @('Egipto', 58141, 0, 15535, 40241, 0, 0, 2365),
    @('Ecuador', 51643, 0, 24991, 22378, 0, 0, 4274),
    @('Paises Bajos', 49804, 82, 0, 0, 0, 2, 6097),
    @('Indonesia', 49009, 1113, 19658, 26778, 0, 38, 2573),
    @('Argentina', 47216, 13, 13816, 32315, 0, 7, 1085),
    @('Emiratos Arabes Unidos', 46133, 450, 34405, 11421, 0, 2, 307),
    @('Singapur', 42623, 191, 36299, 6298, 0, 0, 26),
    @('Kuwait', 41879, 846, 32809, 8733, 0, 3, 337),
    @('Portugal', 40104, 367, 26083, 12478, 0, 3, 1543),
    @('Ucrania', 39014, 940, 17409, 20554, 0, 16, 1051),
    @('Irak', 36702, 2200, 16814, 18558, 0, 79, 1330),
    @('Oman', 33536, 1142, 17972, 15422, 0, 2, 142),
    @('Polonia', 32821, 294, 18134, 13291, 0, 21, 1396),
    @('Filipinas', 32295, 470, 8656, 22435, 0, 18, 1204),
    @('Suiza', 31376, 44, 29000, 418, 0, 2, 1958),
    @('Afganistan', 29640, 159, 9869, 19132, 0, 21, 639),
    @('Republica Dominicana', 28631, 695, 16006, 11934, 0, 16, 691),
    @('Panama', 27314, 0, 14694, 12084, 0, 0, 536),
    @('Bolivia', 26389, 896, 6300, 19243, 0, 26, 846),
    @('Irlanda', 25396, 5, 23364, 306, 0, 6, 1726),
    @('Rumania', 24826, 321, 17391, 5880, 0, 16, 1555),
    @('Barein', 23062, 0, 17450, 5544, 0, 1, 68),
    @('Israel', 21732, 220, 15904, 5520, 0, 0, 308),
    @('Armenia', 21717, 711, 10797, 10534, 0, 14, 386),
    @('Nigeria', 21371, 0, 7338, 13500, 0, 0, 533),
    @('Kazajistan', 18765, 534, 11882, 6747, 0, 2, 136),
    @('Japon', 17968, 0, 16212, 801, 0, 0, 955),
    @('Austria', 17449, 41, 16282, 474, 0, 0, 693),
    @('Moldavia', 15078, 364, 8400, 6183, 0, 5, 495),
    @('Ghana', 14568, 0, 10907, 3566, 0, 0, 95),
    @('Guatemala', 14540, 771, 2897, 11061, 0, 35, 582),
    @('Azerbaiyan', 14305, 590, 7768, 6363, 0, 7, 174),
    @('Honduras', 13943, 587, 1461, 12077, 0, 10, 405),
    @('Serbia', 13235, 143, 12111, 861, 0, 0, 263),
    @('Dinamarca', 12615, 54, 11422, 590, 0, 0, 603),
    @('Corea del Sur', 12535, 51, 10930, 1324, 0, 0, 281),
    @('Camerun', 12270, 0, 7774, 4183, 0, 0, 313),
    @('Argelia', 12248, 172, 8792, 2587, 0, 8, 869),
    @('Marruecos', 10907, 563, 8468, 2223, 0, 2, 216),
    @('Nepal', 10728, 629, 2338, 8366, 0, 0, 24),
    @('Chequia', 10728, 78, 7588, 2797, 0, 4, 343),
    @('Sudan', 8889, 0, 3699, 4642, 0, 0, 548),
    @('Noruega', 8788, 16, 8138, 401, 0, 1, 249),
    @('Malasia', 8596, 6, 8231, 244, 0, 0, 121),
    @('Costa de Marfil', 7904, 0, 3182, 4664, 0, 0, 58),
    @('Australia', 7521, 29, 6915, 503, 0, 1, 103),
    @('Finlandia', 7167, 12, 6600, 240, 0, 0, 327),
    @('Uzbekistan', 6847, 185, 4685, 2143, 0, 0, 19),
    @('Consejo Danes para los Refugiados', 6213, 186, 870, 5201, 0, 7, 142),
    @('Senegal', 6129, 95, 4072, 1964, 0, 4, 93),
    @('Tayikistan', 5630, 63, 4194, 1384, 0, 0, 52),
    @('Republica de Macedonia', 5445, 134, 2091, 3095, 0, 8, 259),
    @('Haiti', 5324, 113, 436, 4799, 0, 1, 89),
    @('Kenia', 5206, 254, 1823, 3253, 0, 2, 130),
    @('El Salvador', 5150, 177, 2924, 2107, 0, 6, 119),
    @('Guinea', 5104, 64, 3744, 1331, 0, 1, 29),
    @('Etiopia', 5034, 186, 1486, 3470, 0, 3, 78),
    @('Gabon', 4849, 0, 2107, 2703, 0, 0, 39),
    @('Republica de Yibuti', 4630, 13, 4182, 396, 0, 3, 52),
    @('Venezuela', 4187, 0, 1327, 2825, 0, 0, 35),
    @('Luxemburgo', 4140, 7, 3965, 65, 0, 0, 110),
    @('Bulgaria', 4114, 0, 2217, 1689, 0, 0, 208),
    @('Hungria', 4114, 7, 2618, 920, 0, 3, 576),
    @('Kirguistan', 3726, 207, 2082, 1602, 0, 1, 42),
    @('Bosnia y Herzegovina', 3676, 88, 2297, 1206, 0, 1, 173),
    @('Grecia', 3310, 8, 1374, 1746, 0, 0, 190),
    @('Mauritania', 3292, 0, 963, 2215, 0, 0, 114),
    @('Tailandia', 3157, 1, 3026, 73, 0, 0, 58),
    @('Republica de Africa Central', 3051, 0, 522, 2492, 0, 0, 37),
    @('Somalia', 2835, 23, 829, 1916, 0, 0, 90),
    @('Guayana Francesa', 2593, 0, 1016, 1569, 0, 0, 8),
    @('Mayotte', 2434, 0, 2188, 215, 0, 0, 31),
    @('Croacia', 2388, 22, 2145, 136, 0, 0, 107),
    @('Costa Rica', 2368, 0, 1129, 1227, 0, 0, 12),
    @('Cuba', 2319, 1, 2130, 104, 0, 0, 85),
    @('Maldivas', 2261, 23, 1839, 414, 0, 0, 8),
    @('Nicaragua', 2170, 0, 1238, 858, 0, 0, 74),
    @('Albania', 2114, 67, 1217, 850, 0, 2, 47),
    @('Mali', 2001, 23, 1333, 556, 0, 1, 112),
    @('Sri Lanka', 2001, 10, 1562, 428, 0, 0, 11),
    @('Estonia', 1983, 1, 1783, 131, 0, 0, 69),
    @('Sudan del Sur', 1930, 0, 217, 1677, 0, 0, 36),
    @('Islandia', 1824, 0, 1806, 8, 0, 0, 10),
    @('Lituania', 1804, 1, 1484, 242, 0, 1, 78),
    @('Madagascar', 1787, 63, 779, 992, 0, 1, 16),
    @('Guinea Ecuatorial', 1664, 0, 515, 1117, 0, 0, 32),
    @('Libano', 1644, 22, 1103, 508, 0, 1, 33),
    @('Eslovaquia', 1607, 18, 1448, 131, 0, 0, 28),
    @('Guinea-Bisau', 1556, 0, 191, 1346, 0, 0, 19),
    @('Eslovenia', 1541, 7, 1376, 56, 0, 0, 109),
    @('Nueva Zelanda', 1516, 1, 1483, 11, 0, 0, 22),
    @('Zambia', 1489, 12, 1223, 248, 0, 0, 18),
    @('Paraguay', 1422, 0, 926, 483, 0, 0, 13),
    @('Sierra Leona', 1354, 7, 869, 430, 0, 0, 55),
    @('Estado de Palestina', 1328, 159, 442, 883, 0, 0, 3),
    @('Hong Kong', 1180, 2, 1086, 88, 0, 0, 6),
    @('Tunez', 1160, 1, 1023, 87, 0, 0, 50),
    @('Letonia', 1111, 0, 903, 178, 0, 0, 30),
    @('Congo', 1087, 0, 456, 594, 0, 0, 37),
    @('Jordania', 1071, 24, 782, 280, 0, 0, 9),
    @('Niger', 1051, 0, 913, 71, 0, 0, 67),
    @('Yemen', 1015, 23, 379, 362, 0, 13, 274),
    @('Cabo Verde', 999, 17, 479, 512, 0, 0, 8),
    @('Republica de Chipre', 990, 0, 824, 147, 0, 0, 19),
    @('Burkina Faso', 919, 12, 825, 41, 0, 0, 53),
    @('Georgia', 914, 3, 771, 129, 0, 0, 14),
    @('Benin', 902, 52, 277, 612, 0, 0, 13),
    @('Uruguay', 885, 0, 815, 45, 0, 0, 25),
    @('Republica del Chad', 860, 0, 757, 29, 0, 0, 74),
    @('Principado de Andorra', 855, 0, 797, 6, 0, 0, 52),
    @('Malaui', 848, 45, 259, 578, 0, 0, 11),
    @('Uganda', 805, 8, 717, 88, 0, 0, 0),
    @('Ruanda', 798, 0, 371, 425, 0, 0, 2),
    @('Mozambique', 762, 5, 220, 537, 0, 0, 5),
    @('Crucero', 712, 0, 651, 48, 0, 0, 13),
    @('Santo Tome y Principe', 707, 0, 208, 487, 0, 0, 12),
    @('San Marino', 698, 0, 647, 9, 0, 0, 42),
    @('Suazilandia', 674, 0, 319, 348, 0, 0, 7),
    @('Jamaica', 670, 5, 518, 142, 0, 0, 10),
    @('Malta', 665, 0, 624, 32, 0, 0, 9),
    @('Liberia', 652, 0, 270, 348, 0, 0, 34),
    @('Libia', 639, 0, 132, 490, 0, 0, 17),
    @('Togo', 576, 0, 384, 179, 0, 0, 13),
    @('Zimbabue', 525, 0, 64, 455, 0, 0, 6),
    @('Tanzania', 509, 0, 183, 305, 0, 0, 21),
    @('Reunion', 507, 0, 460, 46, 0, 0, 1),
    @('Taiwan', 446, 0, 435, 4, 0, 0, 7),
    @('Montenegro', 389, 11, 315, 65, 0, 0, 9),
    @('Vietnam', 352, 3, 329, 23, 0, 0, 0),
    @('Mauricio', 341, 1, 326, 5, 0, 0, 10),
    @('Isla de Man', 336, 0, 312, 0, 0, 0, 24),
    @('Surinam', 319, 0, 132, 178, 0, 0, 9),
    @('Birmania', 293, 1, 208, 79, 0, 0, 6),
    @('Comoras', 265, 0, 159, 99, 0, 0, 7),
    @('Martinica', 236, 0, 98, 124, 0, 0, 14),
    @('Siria', 231, 0, 94, 130, 0, 0, 7),
    @('Mongolia', 215, 0, 168, 47, 0, 0, 0),
    @('Guyana', 206, 0, 107, 87, 0, 0, 12),
    @('Islas Caimanes', 195, 0, 154, 40, 0, 0, 1),
    @('Angola', 189, 0, 77, 102, 0, 0, 10),
    @('Islas Feroe', 187, 0, 187, 0, 0, 0, 0),
    @('Gibraltar', 176, 0, 176, 0, 0, 0, 0),
    @('Guadalupe', 174, 0, 157, 3, 0, 0, 14),
    @('Bermudas', 146, 0, 132, 5, 0, 0, 9),
    @('Burundi', 144, 0, 93, 50, 0, 0, 1),
    @('Eritrea', 143, 0, 39, 104, 0, 0, 0),
    @('Brunei', 141, 0, 138, 0, 0, 0, 3),
    @('Camboya', 130, 0, 127, 3, 0, 0, 0),
    @('Trinidad yTobago', 123, 0, 109, 6, 0, 0, 8),
    @('Bahamas', 104, 0, 83, 10, 0, 0, 11),
    @('Monaco', 102, 1, 95, 3, 0, 0, 4),
    @('Aruba', 101, 0, 98, 0, 0, 0, 3),
    @('Barbados', 97, 0, 85, 5, 0, 0, 7),
    @('Botsuana', 89, 0, 25, 63, 0, 0, 1),
    @('Liechtenstein', 82, 0, 55, 26, 0, 0, 1),
    @('San Martin (Parte Holandesa)', 77, 0, 62, 0, 0, 0, 15),
    @('Namibia', 76, 4, 21, 55, 0, 0, 0),
    @('Butan', 70, 1, 34, 36, 0, 0, 0),
    @('Antigua y Barbuda', 65, 39, 22, 40, 0, 0, 3),
    @('Polinesia Francesa', 60, 0, 60, 0, 0, 0, 0),
    @('Macao', 45, 0, 45, 0, 0, 0, 0),
    @('Gambia', 42, 1, 26, 14, 0, 0, 2),
    @('San Martin (Parte Francesa)', 42, 0, 36, 3, 0, 0, 3),
    @('Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
    @('Guam', 32, 0, 0, 31, 0, 0, 1),
    @('San Vicente y las Granadinas', 29, 0, 27, 2, 0, 0, 0),
    @('Timor Oriental', 24, 0, 24, 0, 0, 0, 0),
    @('Belice', 23, 0, 17, 4, 0, 0, 2),
    @('Curazao', 23, 0, 19, 3, 0, 0, 1),
    @('Granada', 23, 0, 23, 0, 0, 0, 0),
    @('Nueva Caledonia', 21, 0, 21, 0, 0, 0, 0),
    @('Santa Lucia', 19, 0, 18, 1, 0, 0, 0),
    @('Laos', 19, 0, 19, 0, 0, 0, 0),
    @('Fiyi', 18, 0, 18, 0, 0, 0, 0),
    @('Dominica', 18, 0, 18, 0, 0, 0, 0),
    @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
    @('Lesoto', 17, 0, 2, 15, 0, 0, 0),
    @('San Cristobal y Nieves', 15, 0, 15, 0, 0, 0, 0),
    @('Islas Turcas y Caicos', 14, 0, 11, 2, 0, 0, 1),
    @('Islas Malvinas', 13, 0, 13, 0, 0, 0, 0),
    @('Groenlandia', 13, 0, 13, 0, 0, 0, 0),
    @('Santa Sede', 12, 0, 12, 0, 0, 0, 0),
    @('Seychelles', 11, 0, 11, 0, 0, 0, 0),
    @('Montserrat', 11, 0, 10, 0, 0, 0, 1),
    @('Sahara Occidental', 10, 0, 8, 1, 0, 0, 1),
    @('Papua Nueva Guinea', 9, 0, 8, 1, 0, 0, 0),
    @('Islas Virgenes Britanicas', 8, 0, 7, 0, 0, 0, 1),
    @('Bonaire, San Eustaquio y Saba', 7, 0, 7, 0, 0, 0, 0),
    @('San Bartolome', 6, 0, 6, 0, 0, 0, 0),
    @('Anguila', 3, 0, 3, 0, 0, 0, 0),
    @('San Pedro y Miquelon', 1, 0, 1, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 4 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
}

Write-Host "Updated" $data.Count "rows"
